$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-generated data table for rows 76-91 (target state after edit)
$rowsData = @(
    @{ Row=76; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44995; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Especial'; M=56; N=18000; O=18000; P=18000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=1000; T=18 }
    @{ Row=77; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44995; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Primera'; M=50; N=16000; O=16000; P=16000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=889; T=18 }
    @{ Row=78; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44995; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Segunda'; M=48; N=14000; O=14000; P=14000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=778; T=18 }
    @{ Row=79; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44286; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Especial'; M=75; N=12000; O=12000; P=12000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=667; T=18 }
    @{ Row=80; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44286; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Extra (doble especial)'; M=65; N=14000; O=14000; P=14000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=778; T=18 }
    @{ Row=81; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44286; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Primera'; M=70; N=10000; O=10000; P=10000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=556; T=18 }
    @{ Row=82; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44655; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Especial'; M=56; N=15000; O=15000; P=15000; Q='$/caja 15 kilos granel'; R='Región de O''Higgins'; S=1000; T=15 }
    @{ Row=83; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44655; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Primera'; M=67; N=13000; O=13000; P=13000; Q='$/caja 15 kilos granel'; R='Región de O''Higgins'; S=867; T=15 }
    @{ Row=84; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44655; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Segunda'; M=60; N=11000; O=11000; P=11000; Q='$/caja 15 kilos granel'; R='Región de O''Higgins'; S=733; T=15 }
    @{ Row=85; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44637; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Especial'; M=75; N=18000; O=18000; P=18000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=1000; T=18 }
    @{ Row=86; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44637; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Primera'; M=70; N=16000; O=16000; P=16000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=889; T=18 }
    @{ Row=87; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44270; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Especial'; M=65; N=14000; O=14000; P=14000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=778; T=18 }
    @{ Row=88; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44270; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Primera'; M=60; N=12000; O=12000; P=12000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=667; T=18 }
    @{ Row=89; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44273; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Especial'; M=50; N=14000; O=14000; P=14000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=778; T=18 }
    @{ Row=90; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44273; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Primera'; M=60; N=12000; O=12000; P=12000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=667; T=18 }
    @{ Row=91; A=3; B='Femacal de La Calera'; C='Coquimbo'; D=44273; E=5; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; L='Segunda'; M=50; N=10000; O=10000; P=10000; Q='$/caja 18 kilos empedrada'; R='Región de O''Higgins'; S=556; T=18 }
)

$colMap = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10
    K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20
}

foreach ($item in $rowsData) {
    $r = $item.Row
    foreach ($col in $colMap.Keys) {
        $c = $colMap[$col]
        $val = $item[$col]
        $ws.Cells.Item($r, $c).Value = $val
    }
    # Ensure the date column keeps the workbook's date/time number format
    # (matches the formatting already used for column D elsewhere in the sheet).
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Done writing rows 76-91"
Write-Host "New dimension:" $ws.UsedRange.Address()
